# Updates the cryptocurrency price/volume tracker sheet with the latest
# scraped figures (GitHub Actions scheduled refresh), including two rows
# whose coin ranking swapped position relative to each other.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume) hold plain-text figures (dotted
# thousands separators, +/-% strings, no real numeric semantics). Force
# text format first so Excel does not reinterpret numeric-looking values
# (e.g. "1.00", "13.50", "0.0000360") as actual numbers and strip formatting.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.219.69'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.81%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.923.60'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.76%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '488.66'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +4.14%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.86'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +2.17%  '

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.34%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.04%  '

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.93%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.171'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +4.21%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0000360'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +7.40%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '42.81'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.22%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.68'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +2.85%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.550.21'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.90%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.91'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.62%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.931.85'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.61%  '

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.06%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '20.13'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.17%  '

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.75%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '68.352.01'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.78%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '446.77'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +3.69%  '

$ws.Range("B22").Value = 'ImmutableX'
$ws.Range("C22").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.41'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +2.43%  '

$ws.Range("B23").Value = 'InternetComputer(DFINITY)'
$ws.Range("C23").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.81'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.66%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '88.55'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.18%  '

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +15.43%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.01'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +14.87%  '

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +3.00%  '

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.63%  '

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +2.10%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '702.24'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -3.76%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '13.50'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.56%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.132'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.59%  '

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +2.58%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0₃0941'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +21.16%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '41.82'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -2.89%  '

$ws.Range("B36").Value = 'OKB'
$ws.Range("C36").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '58.95'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +2.70%  '

$ws.Range("B37").Value = 'NEARProtocol'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.82'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +7.78%  '

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -4.12%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.00'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.02%  '

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.80%  '

$ws.Range("B41").Value = 'Fetch.AI'
$ws.Range("C41").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.87'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +14.00%  '

$ws.Range("B42").Value = 'TheGraph'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.375'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +10.78%  '

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -3.61%  '

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +5.77%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.144'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +2.18%  '

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.04%  '

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.99%  '

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.80%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '146.15'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.70%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.15'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.85%  '

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.54%  '
